$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

$updates = @{
    "H2" = 428501.4
    "J2" = 1360.5
    "L2" = 1360.5
    "N2" = -1586.5
    "H21" = 1402
    "I21" = 1402
    "J21" = 0
    "K21" = 1402
    "L21" = 0
    "M21" = -934
    "H23" = 1402
    "I23" = 1402
    "J23" = 0
    "K23" = 1402
    "L23" = 0
    "M23" = -1168
    "H29" = 666
    "I29" = 666
    "J29" = 0
    "K29" = 1998
    "L29" = 0
    "M29" = -1717
    "H32" = 3424
    "J32" = 3632.3333
    "L32" = 3632.3333
    "N32" = -4284.3333
    "H38" = 372.73334
    "I38" = 85.07143000000001
    "J38" = 4400
    "K38" = 255.21429
    "L38" = 13200
    "M38" = 116.78571
    "N38" = -13944
    "H58" = 1162.1818
    "I58" = 1162.1818
    "K58" = 3486.5454
    "M58" = -3336.5454
    "H62" = 8100
    "I62" = 7093.4614
    "K62" = 7093.4614
    "M62" = -6469.4614
    "H65" = 8100
    "I65" = 7093.4614
    "K65" = 35467.307
    "M65" = -32347.307
    "H76" = 5498.1113
    "I76" = 5748
    "K76" = 5748
    "M76" = -5433
    "H79" = 5498.1113
    "I79" = 5748
    "K79" = 5748
    "M79" = -4656
    "H80" = 1283.3334
    "I80" = 967
    "K80" = 2901
    "M80" = -1903
    "H83" = 1283.3334
    "I83" = 967
    "K83" = 8703
    "M83" = -3711
    "H106" = 5071.3887
    "I106" = 3294.375
    "K106" = 3294.375
    "M106" = -2663.375
    "H108" = 99000
    "J108" = 99000
    "L108" = 99000
    "N108" = -106680
    "H112" = 3315.1428
    "J112" = 3315.1428
    "L112" = 9945.428400000001
    "N112" = -12161.4284
    "H113" = 4775.2856
    "I113" = 4491.467
    "J113" = 5484.8335
    "K113" = 4491.467
    "L113" = 5484.8335
    "M113" = -1237.467
    "N113" = -11992.8335
    "H132" = 3041.0286
    "I132" = 2607.9644
    "J132" = 4773.2856
    "K132" = 7823.8932
    "L132" = 14319.8568
    "M132" = -5293.8932
    "N132" = -19379.8568
    "H135" = 538.4286
    "I135" = 474.33334
    "K135" = 4269.00006
    "M135" = -1734.00006
    "H137" = 26672.75
    "I137" = 42955.043
    "J137" = 2249.3125
    "K137" = 128865.129
    "L137" = 6747.9375
    "M137" = -126315.129
    "N137" = -11847.9375
    "H138" = 1847.6459
    "I138" = 1654.0286
    "J138" = 2368.923
    "K138" = 4962.085800000001
    "L138" = 7106.768999999999
    "M138" = 177.9141999999993
    "N138" = -17386.769
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

foreach ($ref in @("N21", "N23", "N29")) {
    $ws.Range($ref).ClearContents()
}

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

$updates = @{
    "H45" = 4871.2
    "I45" = 1300
    "K45" = 1300
    "M45" = -923
    "H61" = 7947.6855
    "I61" = 9704.888999999999
    "J61" = 2017.125
    "K61" = 9704.888999999999
    "L61" = 2017.125
    "M61" = -9492.888999999999
    "N61" = -2441.125
    "H64" = 200000
    "I64" = 0
    "J64" = 200000
    "K64" = 0
    "L64" = 200000
    "N64" = -200496
    "H67" = 200000
    "I67" = 0
    "J67" = 200000
    "K67" = 0
    "L67" = 200000
    "N67" = -201716
    "H74" = 31776.146
    "I74" = 36861.31
    "K74" = 36861.31
    "M74" = -35987.31
    "H77" = 31776.146
    "I77" = 36861.31
    "K77" = 184306.55
    "M77" = -179938.55
    "H132" = 29179.79
    "I132" = 32973.91
    "J132" = 4138.6
    "K132" = 98921.73000000001
    "L132" = 12415.8
    "M132" = -96391.73000000001
    "N132" = -17475.8
    "H136" = 7947.6855
    "I136" = 9704.888999999999
    "J136" = 2017.125
    "K136" = 29114.667
    "L136" = 6051.375
    "M136" = -26564.667
    "N136" = -11151.375
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

foreach ($ref in @("M64", "M67")) {
    $ws.Range($ref).ClearContents()
}

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

$updates = @{
    "H22" = 59663.65
    "I22" = 91707.73
    "K22" = 91707.73
    "M22" = -91534.73
    "H80" = 665
    "I80" = 731.25
    "K80" = 731.25
    "M80" = 266.75
    "H83" = 665
    "I83" = 731.25
    "K83" = 3656.25
    "M83" = 1335.75
    "H86" = 2452.3845
    "I86" = 1906.05
    "J86" = 4273.5
    "K86" = 1906.05
    "L86" = 4273.5
    "M86" = -783.05
    "N86" = -6519.5
    "H89" = 2452.3845
    "I89" = 1906.05
    "J89" = 4273.5
    "K89" = 9530.25
    "L89" = 21367.5
    "M89" = -3914.25
    "N89" = -32599.5
    "H99" = 2765.0625
    "J99" = 2355.6
    "L99" = 2355.6
    "N99" = -5351.6
    "H105" = 3525
    "I105" = 2993.2856
    "K105" = 2993.2856
    "M105" = -1246.2856
    "H134" = 2640.3684
    "I134" = 2390.5715
    "K134" = 7171.7145
    "M134" = -4636.7145
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

$updates = @{
    "H22" = 384
    "I22" = 255.5
    "J22" = 589.6
    "K22" = 255.5
    "L22" = 589.6
    "M22" = 94.5
    "N22" = -1289.6
    "H31" = 2898.9355
    "I31" = 2409.2415
    "J31" = 9999.5
    "K31" = 2409.2415
    "L31" = 9999.5
    "M31" = -2114.2415
    "N31" = -10589.5
    "H34" = 2898.9355
    "I34" = 2409.2415
    "J34" = 9999.5
    "K34" = 2409.2415
    "L34" = 9999.5
    "M34" = -2207.2415
    "N34" = -10403.5
    "H105" = 41379.4
    "I105" = 41379.4
    "K105" = 41379.4
    "M105" = -39632.4
    "H132" = 1089.9395
    "I132" = 1089.5161
    "K132" = 3268.5483
    "M132" = -738.5483000000004
    "H134" = 35788.656
    "I134" = 38339.85
    "J134" = 1347.5
    "K134" = 115019.55
    "L134" = 4042.5
    "M134" = -112484.55
    "N134" = -9112.5
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

$updates = @{
    "H4" = 1172881.1
    "I4" = 864861.6
    "K4" = 2594584.8
    "M4" = -2594472.8
    "H38" = 244.45454
    "I38" = 185.71428
    "K38" = 557.14284
    "M38" = -210.14284
    "H103" = 1023.6875
    "I103" = 477.33334
    "J103" = 1351.5
    "K103" = 1432.00002
    "L103" = 4054.5
    "M103" = -553.0000199999999
    "N103" = -5812.5
    "H132" = 3128.5715
    "I132" = 1980
    "K132" = 17820
    "M132" = -15290
    "H141" = 1086.75
    "I141" = 1086.75
    "K141" = 3260.25
    "M141" = 1919.75
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

$updates = @{
    "H11" = 2000
    "J11" = 2000
    "L11" = 2000
    "N11" = -2278
    "H26" = 37899.5
    "J26" = 37899.5
    "L26" = 37899.5
    "N26" = -38459.5
    "H50" = 37899.5
    "J50" = 37899.5
    "L50" = 37899.5
    "N50" = -38895.5
    "H102" = 2903
    "I102" = 2903
    "K102" = 2903
    "M102" = -1281
    "H132" = 33885.918
    "I132" = 38864.29
    "K132" = 116592.87
    "M132" = -114062.87
    "H136" = 27477.723
    "J136" = 27477.723
    "L136" = 82433.16900000001
    "N136" = -87533.16900000001
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

$updates = @{
    "H7" = 7257.88
    "I7" = 10108.143
    "K7" = 10108.143
    "M7" = -9996.143
    "H16" = 5255.25
    "I16" = 5576.5625
    "J16" = 3970
    "K16" = 5576.5625
    "L16" = 3970
    "M16" = -5406.5625
    "N16" = -4310
    "H40" = 14448.125
    "I40" = 16819.666
    "K40" = 16819.666
    "M40" = -16683.666
    "H46" = 14753.1
    "I46" = 28371.875
    "J46" = 5673.9165
    "K46" = 28371.875
    "L46" = 5673.9165
    "M46" = -28183.875
    "N46" = -6049.9165
    "H68" = 3912.25
    "I68" = 3299.6
    "J68" = 4933.3335
    "K68" = 3299.6
    "L68" = 4933.3335
    "M68" = -2550.6
    "N68" = -6431.3335
    "H71" = 3912.25
    "I71" = 3299.6
    "J71" = 4933.3335
    "K71" = 16498
    "L71" = 24666.6675
    "M71" = -12754
    "N71" = -32154.6675
    "H88" = 1000000
    "I88" = 0
    "J88" = 1000000
    "K88" = 0
    "L88" = 1000000
    "N88" = -1000856
    "H91" = 1000000
    "I91" = 0
    "J91" = 1000000
    "K91" = 0
    "L91" = 1000000
    "N91" = -1002964
    "H93" = 1370.0682
    "I93" = 1357.9412
    "J93" = 1377.7037
    "K93" = 1357.9412
    "L93" = 1377.7037
    "M93" = -109.9412
    "N93" = -3873.7037
    "H126" = 7257.88
    "I126" = 10108.143
    "K126" = 30324.429
    "M126" = -27854.429
    "H132" = 19673.041
    "I132" = 21967.42
    "K132" = 65902.25999999999
    "M132" = -63372.25999999999
    "H136" = 2629.276
    "I136" = 2361.111
    "K136" = 7083.333
    "M136" = -4533.333
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

foreach ($ref in @("M88", "M91")) {
    $ws.Range($ref).ClearContents()
}

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

$updates = @{
    "H54" = 18000
    "I54" = 18000
    "J54" = 18000
    "K54" = 18000
    "L54" = 18000
    "M54" = -17480
    "N54" = -19040
    "H87" = 0
    "J87" = 0
    "L87" = 0
    "H90" = 0
    "J90" = 0
    "L90" = 0
    "H100" = 971.6667
    "I100" = 719.6
    "K100" = 1439.2
    "M100" = -898.2
    "H107" = 2389.6
    "I107" = 649.6667
    "K107" = 1949.0001
    "M107" = -29.00009999999997
    "H132" = 24778.064
    "I132" = 25262.045
    "K132" = 75786.13499999999
    "M132" = -73256.13499999999
    "H136" = 3827.7097
    "I136" = 3275
    "K136" = 9825
    "M136" = -7275
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

foreach ($ref in @("N87", "N90")) {
    $ws.Range($ref).ClearContents()
}
